$wb = $excel.ActiveWorkbook

# --- "forms" sheet: rename the address form's ident from "address-form" to "address" ---
$wsForms = $wb.Worksheets.Item("forms")
$wsForms.Range("A2").Value = 'address'
$wsForms.Range("B2").Value = '{"ident":"address","elems":[{"ident":"TITLE-60664","elementType":"TITLE","texts":{"label":{"textType":"LABEL","texts":{"EN":"Address","DE":"Adresse"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE_CLASS","value":"huge"}]},"value":"","required":false,"inline":false,"readOnly":false,"layoutWide":"SIXTEEN","elemEntries":{"entries":[]},"validations":{"rules":[]}},{"ident":"street","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Street","DE":"Strasse"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE","value":"20"},{"extraProp":"INPUT_TYPE","value":"text"}]},"value":"Sonnenweg","required":true,"inline":false,"readOnly":false,"layoutWide":"TWELVE","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"number","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Number","DE":"Nummer"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE","value":"20"},{"extraProp":"INPUT_TYPE","value":"text"}]},"value":"23a","required":true,"inline":false,"readOnly":false,"layoutWide":"FOUR","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"postcode","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Postcode","DE":"Plz"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE","value":"20"},{"extraProp":"INPUT_TYPE","value":"number"}]},"value":"6414","required":true,"inline":false,"readOnly":false,"layoutWide":"FOUR","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":true,"params":{"intParam1":1000,"intParam2":9999}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"town","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Town","DE":"Ort"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE","value":"20"},{"extraProp":"INPUT_TYPE","value":"text"}]},"value":"Oberarth","required":true,"inline":false,"readOnly":false,"layoutWide":"TWELVE","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}}]}'

# --- "data" sheet: restructure address-data's JSON (object -> ordered pairs, city -> town) ---
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("B2").Value = '{"ident":"address-data","structure":{"value":[["street",{"DataString":{"value":"Sonnenweg"}}],["number",{"DataString":{"value":"23a"}}],["postcode",{"DataNumber":{"value":6414}}],["town",{"DataString":{"value":"Oberarth"}}]]}}'

# --- "mappings" sheet: auto-derive mapping idents directly from form/data field names ---
$wsMappings = $wb.Worksheets.Item("mappings")
$wsMappings.Range("B2").Value = '{"ident":"address-mapping","formIdent":"address","dataIdent":"address-data","mappings":[{"formIdent":"street","dataIdent":"street"},{"formIdent":"number","dataIdent":"number"},{"formIdent":"postcode","dataIdent":"postcode"},{"formIdent":"town","dataIdent":"town"}]}'

# --- Update the view/selection state on each affected sheet ---
$wsForms.Select()
$wsForms.Range("A6").Select()

$wsData.Select()
$wsData.Range("B2").Select()

$wsMappings.Select()
$wsMappings.Range("B6").Select()
